# Update the "Förändrad" (changed) date column (C) for rows 2-5
# from serial date 45233 (2023-11-03) to 45243 (2023-11-13).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($row in 2..5) {
    $ws.Cells.Item($row, 3).Value = 45243
}
